# 713 inserted MY record
#
# A new "MY" results column is inserted between column A (LOC bucket) and
# the old "SourcerCC" column of the small results table that lives in
# A14:E20. Everything that used to sit in B:D shifts one column to the
# right (B->C, C->D, D->E) and the new "MY" timings are written into the
# now-empty column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timing values for the inserted "MY" column, keyed by row.
$myValues = @{
    14 = "MY"
    15 = "1s"
    16 = "1s"
    17 = "4s"
    18 = "9s"
    19 = "59s"
    20 = "12min23s"
}

for ($row = 14; $row -le 20; $row++) {
    # Shift existing D -> E, C -> D, B -> C (right to left so we never
    # clobber a value before it has been copied onward). Reads must go
    # through .Value2 - plain .Value getter is not reliable here.
    $ws.Cells.Item($row, 5).Value = $ws.Cells.Item($row, 4).Value2
    $ws.Cells.Item($row, 4).Value = $ws.Cells.Item($row, 3).Value2
    $ws.Cells.Item($row, 3).Value = $ws.Cells.Item($row, 2).Value2

    # Write the new "MY" value into the freed-up column B.
    $ws.Cells.Item($row, 2).Value = $myValues[$row]
}

# Row 14 is the header row (A14/C14/D14/E14 already carry the yellow
# highlight style) - give the new B14 header cell the same fill so it
# matches its neighbours.
$ws.Range("B14").Interior.Color = 65535

# The active selection moved from the old header cell E20 to the newly
# inserted B20 while editing the table.
[void]$ws.Range("B20").Select()
